$wb = $excel.ActiveWorkbook

# ----- Sheet "daily": add the 8 new vocabulary rows (A:B), rows 2-9 -----
$wsDaily = $wb.Worksheets.Item("daily")

$wsDaily.Range("A2").Value = "se lever"
$wsDaily.Range("B2").Value = "get up"
$wsDaily.Range("A3").Value = "s'endormir"
$wsDaily.Range("B3").Value = "fall asleep"
$wsDaily.Range("A4").Value = "dormir"
$wsDaily.Range("B4").Value = "sleep"
$wsDaily.Range("A5").Value = "se réveiller"
$wsDaily.Range("B5").Value = "wake up"
$wsDaily.Range("A6").Value = "manger"
$wsDaily.Range("B6").Value = "eat"
$wsDaily.Range("A7").Value = "boire"
$wsDaily.Range("B7").Value = "drink"
$wsDaily.Range("A8").Value = "coucher"
$wsDaily.Range("B8").Value = "go bed"
$wsDaily.Range("A9").Value = "nettoyer"
$wsDaily.Range("B9").Value = "clean"

# match the data-row look (centered horizontal + vertical) used elsewhere in the workbook
$wsDaily.Range("A2:B9").HorizontalAlignment = -4108
$wsDaily.Range("A2:B9").VerticalAlignment = -4108

# column widths for A/B on the "daily" sheet
$wsDaily.Columns.Item(1).ColumnWidth = 10.428571428571429
$wsDaily.Columns.Item(2).ColumnWidth = 12.571428571428573

# ----- Sheet "mouth": the manger/eat + boire/drink rows moved to "daily", clear them here -----
$wsMouth = $wb.Worksheets.Item("mouth")
$wsMouth.Range("A12:B13").ClearContents()

# ----- view/selection updates -----
$wsMouth.Activate()
$excel.ActiveWindow.Zoom = 180
$wsMouth.Range("A12:B13").Select()

$wsDaily.Activate()
$excel.ActiveWindow.Zoom = 227
$wsDaily.Range("B6").Select()
